# buffalocart error msg on resetpassword test
#
# 1) rename sheet "home" -> "Recoveremail" and make it the active tab
# 2) Login sheet: rename header cells ("valid username"/"valid password" -> "username"/"password")
#    and re-point the selection to C5
# 3) Recoveremail sheet: fill in recover-email-id / reset-error-msg data, add a
#    mailto hyperlink on the e-mail cell, and re-point the selection to C7

$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item(1)
$wsRecover = $wb.Worksheets.Item(2)

# --- Login sheet -----------------------------------------------------------
$wsLogin.Range("B1").Value = "username"
$wsLogin.Range("C1").Value = "password"

$wsLogin.Range("A2").Value = "Login - Demo POS"
$wsLogin.Range("B2").Value = "admin"
$wsLogin.Range("D2").Value = "Aju Mathew"
$wsLogin.Range("E2").Value = "These credentials do not match our records."

[void]$wsLogin.Range("C5").Select()

# --- Recoveremail sheet -----------------------------------------------------
$wsRecover.Name = "Recoveremail"

$wsRecover.Range("A1").Value = "recover email id"
$wsRecover.Range("B1").Value = "reset error msg"
$wsRecover.Range("A2").Value = "ann76@gmail.com"
$wsRecover.Range("B2").Value = "We can't find a user with that e-mail address."

[void]$wsRecover.Hyperlinks.Add($wsRecover.Range("A2"), "mailto:ann76@gmail.com")

[void]$wsRecover.Range("C7").Select()

# Recoveremail becomes the active/visible tab
[void]$wsRecover.Activate()
